$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The average-score column (B) had its decimal values replaced with plain
# text using a comma as decimal separator (e.g. 7.5 -> "7,5"), and lost the
# number-format style that used to be applied to those cells. Whole number
# rows (B3, B7) are left untouched.
$ws.Range("B2").Value = "7,5"
$ws.Range("B2").ClearFormats()

$ws.Range("B4").Value = "5,25"
$ws.Range("B4").ClearFormats()

$ws.Range("B5").Value = "5,5"
$ws.Range("B5").ClearFormats()

$ws.Range("B6").Value = "7,5"
$ws.Range("B6").ClearFormats()

$ws.Range("B8").Value = "7,25"
$ws.Range("B8").ClearFormats()

$ws.Range("B9").Value = "7,25"
$ws.Range("B9").ClearFormats()
